$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A values to be serial numbers 1..13 instead of the
# duplicated company-name text.
for ($r = 1; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
}

# Column A becomes a narrow index column (~3.14 chars wide); column B is
# left untouched so it keeps the original wide "bestFit" width that used
# to be shared across the combined A:B range.
$ws.Columns.Item(1).ColumnWidth = 2.3

# Update the active selection to the new index column range.
$ws.Range("A1:A13").Select()
